$wb = $excel.ActiveWorkbook

# Update actual/start price values on both the "New" and "LastDownload" sheets
$sheetNames = @("New", "LastDownload")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D40").Value = "621 000,00 ₽"
    $ws.Range("D132").Value = "216 000,00 ₽"
    $ws.Range("D133").Value = "4 500 000,00 ₽"
    $ws.Range("D134").Value = "747 900,00 ₽"
    $ws.Range("D135").Value = "385 200,00 ₽"
    $ws.Range("D150").Value = "2 072 250,00 ₽"
    $ws.Range("D157").Value = "217 800,00 ₽"
    $ws.Range("D158").Value = "270 000,00 ₽"
    $ws.Range("D159").Value = "63 000,00 ₽"
    $ws.Range("D160").Value = "94 500,00 ₽"
    $ws.Range("D179").Value = "131 184,00 ₽"
    $ws.Range("D180").Value = "45 000,00 ₽"
    $ws.Range("D181").Value = "166 353 568,43 ₽"
    $ws.Range("D182").Value = "63 000,00 ₽"
    $ws.Range("D183").Value = "95 400,00 ₽"
    $ws.Range("D184").Value = "213 411,00 ₽"
    $ws.Range("D185").Value = "54 630,00 ₽"
    $ws.Range("D186").Value = "116 091,90 ₽"
    $ws.Range("D187").Value = "1 342 859,49 ₽"
    $ws.Range("D188").Value = "385 560,00 ₽"
    $ws.Range("D189").Value = "3 607 200,00 ₽"
}

# Update the selection on the active sheet ("LastDownload") to reflect the saved cursor position
$wsLastDownload = $wb.Worksheets.Item("LastDownload")
$wsLastDownload.Activate()
$wsLastDownload.Range("M8").Select()

